# Applies the Codebook.xlsx commit: adds the DOMCREDIT-gap / venture-capital /
# skilled-employees codebook rows and restyles the new footnote cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be created in this order so they land on the
# same shared-string indices (26-29) the authored workbook uses.

# F6 / F7: new plain source notes
$ws.Range("F6").Value = "Availability of venture capital indicates how easy it is for entrepreneurs with innovative but risky projects to find venture capital."
$ws.Range("F7").Value = "DOMCREDIT gap"

# F8: new note, rendered in a small Arial footnote font
$f8 = $ws.Range("F8")
$f8.Value = "Financial resources provided to the private sector by financial corporations as a percentage of GDP. [i]Financial resources[i] are loans, purchases of non-equity securities, and trade credits and other accounts receivable, that establish a claim for repayment."
$f8.Font.Name = "Arial"
$f8.Font.Size = 8

# D5: new codebook label for EOSQ403 (ease of finding skilled employees)
$ws.Range("D5").Value = "EOSQ403: Ease of finding skilled employees"

# Row 5 grows a touch taller to fit the new wrapped text
$ws.Range("A5:I5").RowHeight = 28

# Restore the selection anchor left by the author after the edit
$ws.Range("E8").Select()
